# Update "Förändrad" (changed) date column C for all data rows (2-15)
# from serial date 45233 to serial date 45243, leaving everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45233) {
        $cell.Value = 45243
    }
}
